$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12 (shifts old rows 12-26 down to 13-27)
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with values (mirrors the neighboring rows' static
# columns, with its own Fecha/Volumen)
$ws.Cells.Item(12, 1).Value = 8
$ws.Cells.Item(12, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(12, 3).Value = "Coquimbo"
$ws.Cells.Item(12, 4).Value = 45062
$ws.Cells.Item(12, 4).NumberFormat = $ws.Cells.Item(13, 4).NumberFormat
$ws.Cells.Item(12, 5).Value = 4
$ws.Cells.Item(12, 6).Value = 100112039
$ws.Cells.Item(12, 7).Value = "Ciboulette"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 1100
$ws.Cells.Item(12, 11).Value = 2000
$ws.Cells.Item(12, 12).Value = 2500
$ws.Cells.Item(12, 13).Value = 2250
$ws.Cells.Item(12, 14).Value = "$/docena de atados"
$ws.Cells.Item(12, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(12, 16).Value = 750
$ws.Cells.Item(12, 17).Value = 3
$ws.Cells.Item(12, 18).Value = "Hortaliza"
